# Applies the "three-digit number multiplied by one-digit number" update:
# - bumps the header date by one day
# - swaps in a fresh set of multiplication problems throughout the table

$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-18 Tuesday", "2025-02-19 Wednesday"),
    @("473×9=", "315×4="),
    @("944×7=", "371×8="),
    @("251×7=", "147×2="),
    @("556×6=", "133×6="),
    @("432×5=", "700×8="),
    @("283×8=", "776×7="),
    @("601×7=", "712×4="),
    @("418×8=", "980×8="),
    @("788×5=", "541×6="),
    @("638×2=", "719×3="),
    @("511×7=", "576×3="),
    @("353×7=", "142×8="),
    @("323×3=", "122×5="),
    @("598×3=", "604×9="),
    @("647×8=", "263×2="),
    @("531×5=", "344×3="),
    @("438×5=", "304×4="),
    @("146×6=", "832×4="),
    @("366×4=", "132×7="),
    @("126×2=", "745×9="),
    @("589×3=", "882×3="),
    @("597×9=", "388×4="),
    @("267×8=", "147×3="),
    @("432×7=", "104×8="),
    @("582×3=", "459×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
